# Update profit/price figures across several Leve sheets following a
# scheduled data refresh (currentAveragePrice* / LevePrice* / LeveProfit*
# columns, i.e. columns H-N of the affected rows).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 3108.5
$ws.Range("I29").Value = 478.33334
$ws.Range("J29").Value = 10999
$ws.Range("K29").Value = 1435.00002
$ws.Range("L29").Value = 32997
$ws.Range("M29").Value = -1154.00002
$ws.Range("N29").Value = -33559

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# Row 92
$ws.Range("H92").Value = 59475.176
$ws.Range("I92").Value = 77310.766
$ws.Range("J92").Value = 1509.5
$ws.Range("K92").Value = 77310.766
$ws.Range("L92").Value = 1509.5
$ws.Range("M92").Value = -76062.766
$ws.Range("N92").Value = -4005.5

# Row 106
$ws.Range("H106").Value = 5249.7
$ws.Range("I106").Value = 3738.75
$ws.Range("J106").Value = 11293.5
$ws.Range("K106").Value = 3738.75
$ws.Range("L106").Value = 11293.5
$ws.Range("M106").Value = -3107.75
$ws.Range("N106").Value = -12555.5

# Row 107
$ws.Range("H107").Value = 1865.238
$ws.Range("I107").Value = 1832.8334
$ws.Range("J107").Value = 2059.6667
$ws.Range("K107").Value = 1832.8334
$ws.Range("L107").Value = 2059.6667
$ws.Range("M107").Value = 87.16660000000002
$ws.Range("N107").Value = -5899.6667

# Row 137
$ws.Range("H137").Value = 3579.4583
$ws.Range("I137").Value = 3213.923
$ws.Range("J137").Value = 4011.4546
$ws.Range("K137").Value = 9641.769
$ws.Range("L137").Value = 12034.3638
$ws.Range("M137").Value = -7091.769
$ws.Range("N137").Value = -17134.3638

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 14663.191
$ws.Range("I32").Value = 14158.762
$ws.Range("J32").Value = 23137.6
$ws.Range("K32").Value = 14158.762
$ws.Range("L32").Value = 23137.6
$ws.Range("M32").Value = -13871.762
$ws.Range("N32").Value = -23711.6

# Row 61
$ws.Range("H61").Value = 10883.615
$ws.Range("I61").Value = 11540.583
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 11540.583
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -11328.583
$ws.Range("N61").Value = -3424

# Row 132
$ws.Range("H132").Value = 28492.838
$ws.Range("I132").Value = 29978.143
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 89934.429
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -87404.429
$ws.Range("N132").Value = -12560

# Row 136
$ws.Range("H136").Value = 10883.615
$ws.Range("I136").Value = 11540.583
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 34621.749
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -32071.749
$ws.Range("N136").Value = -14100

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2378.8
$ws.Range("I20").Value = 1899.7778
$ws.Range("J20").Value = 3610.5715
$ws.Range("K20").Value = 1899.7778
$ws.Range("L20").Value = 3610.5715
$ws.Range("M20").Value = -1652.7778
$ws.Range("N20").Value = -4104.5715

# Row 134
$ws.Range("H134").Value = 1778.34
$ws.Range("I134").Value = 1687.0889
$ws.Range("J134").Value = 2599.6
$ws.Range("K134").Value = 5061.2667
$ws.Range("L134").Value = 7798.799999999999
$ws.Range("M134").Value = -2526.2667
$ws.Range("N134").Value = -12868.8

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 40000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 40000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 40000
$ws.Range("N41").Value = -40856

# Row 43
$ws.Range("H43").Value = 31332.666
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 31332.666
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 31332.666
$ws.Range("N43").Value = -31700.666

# Row 99
$ws.Range("H99").Value = 6861.316
$ws.Range("I99").Value = 10113.909
$ws.Range("J99").Value = 2389
$ws.Range("K99").Value = 10113.909
$ws.Range("L99").Value = 2389
$ws.Range("M99").Value = -8615.909
$ws.Range("N99").Value = -5385

# Row 101
$ws.Range("H101").Value = 31332.666
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 31332.666
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 31332.666
$ws.Range("N101").Value = -37822.666

# Row 126
$ws.Range("H126").Value = 6861.316
$ws.Range("I126").Value = 10113.909
$ws.Range("J126").Value = 2389
$ws.Range("K126").Value = 30341.727
$ws.Range("L126").Value = 7167
$ws.Range("M126").Value = -27871.727
$ws.Range("N126").Value = -12107

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 946.6
$ws.Range("I8").Value = 946.6
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 2839.8
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2700.8

# Row 132
$ws.Range("H132").Value = 994.5
$ws.Range("I132").Value = 992.6667
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 8934.0003
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6404.0003
$ws.Range("N132").Value = -14060

# Row 133
$ws.Range("H133").Value = 7733.3335
$ws.Range("I133").Value = 5625
$ws.Range("J133").Value = 11950
$ws.Range("K133").Value = 16875
$ws.Range("L133").Value = 35850
$ws.Range("M133").Value = -11815
$ws.Range("N133").Value = -45970

# Row 134
$ws.Range("H134").Value = 1909.6154
$ws.Range("I134").Value = 1626.5454
$ws.Range("J134").Value = 3466.5
$ws.Range("K134").Value = 4879.6362
$ws.Range("L134").Value = 10399.5
$ws.Range("M134").Value = 190.3638000000001
$ws.Range("N134").Value = -20539.5

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# Row 139
$ws.Range("H139").Value = 913.2
$ws.Range("I139").Value = 913.2
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 2739.6
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 2400.4
$ws.Range("N139").ClearContents()

# Row 141
$ws.Range("H141").Value = 1539.75
$ws.Range("I141").Value = 1539.75
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 560.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3335.375
$ws.Range("I80").Value = 2772.25
$ws.Range("J80").Value = 3898.5
$ws.Range("K80").Value = 2772.25
$ws.Range("L80").Value = 3898.5
$ws.Range("M80").Value = -1774.25
$ws.Range("N80").Value = -5894.5

# Row 83
$ws.Range("H83").Value = 3335.375
$ws.Range("I83").Value = 2772.25
$ws.Range("J83").Value = 3898.5
$ws.Range("K83").Value = 13861.25
$ws.Range("L83").Value = 19492.5
$ws.Range("M83").Value = -8869.25
$ws.Range("N83").Value = -29476.5

# Row 131
$ws.Range("H131").Value = 34199.4
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 34199.4
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 34199.4
$ws.Range("N131").Value = -44279.4

# Row 132
$ws.Range("H132").Value = 15947.761
$ws.Range("I132").Value = 17597.746
$ws.Range("J132").Value = 2954.125
$ws.Range("K132").Value = 52793.238
$ws.Range("L132").Value = 8862.375
$ws.Range("M132").Value = -50263.238
$ws.Range("N132").Value = -13922.375

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 35843.64
$ws.Range("I132").Value = 45024.215
$ws.Range("J132").Value = 3711.625
$ws.Range("K132").Value = 135072.645
$ws.Range("L132").Value = 11134.875
$ws.Range("M132").Value = -132542.645
$ws.Range("N132").Value = -16194.875

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 27500.78
$ws.Range("I122").Value = 1416.3662
$ws.Range("J122").Value = 336166.34
$ws.Range("K122").Value = 4249.098599999999
$ws.Range("L122").Value = 1008499.02
$ws.Range("M122").Value = -1799.098599999999
$ws.Range("N122").Value = -1013399.02

# Row 132
$ws.Range("H132").Value = 16608.318
$ws.Range("I132").Value = 20048.822
$ws.Range("J132").Value = 1787.6923
$ws.Range("K132").Value = 60146.466
$ws.Range("L132").Value = 5363.0769
$ws.Range("M132").Value = -57616.466
$ws.Range("N132").Value = -10423.0769

# Row 136
$ws.Range("H136").Value = 2599.52
$ws.Range("I136").Value = 2414.2632
$ws.Range("J136").Value = 3186.1667
$ws.Range("K136").Value = 7242.7896
$ws.Range("L136").Value = 9558.500100000001
$ws.Range("M136").Value = -4692.7896
$ws.Range("N136").Value = -14658.5001
